# Add "stakoverflow" and "wiki" hyperlinked reference textboxes to slide 16,
# and resize/reposition the background picture on that slide.
#
# EMU -> point conversion helper (PowerPoint Shape geometry is expressed in
# points over COM, while the OOXML stores EMU; 1 pt = 12700 EMU).
function EMU([double]$v) { return $v / 12700.0 }

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)

# --- resize/reposition the full-bleed background picture (Picture 2) -------
$bg = $s.Shapes.Item(1)
$bg.Left   = (EMU 0)
$bg.Top    = (EMU -1)
$bg.Width  = (EMU 12059728)
$bg.Height = (EMU 5872565)

# --- "stakoverflow" textbox (links to Stack Overflow) -----------------------
$tb1 = $s.Shapes.AddTextbox(1, (EMU 9420045), (EMU 5171288), (EMU 1982017), (EMU 400110))
$tb1.Name = "TextBox 1"
$tb1.Fill.Visible = 0
$tb1.TextFrame.WordWrap = 0
$tb1.TextFrame.AutoSize = 1

$tr1 = $tb1.TextFrame.TextRange
$tr1.Text = "stakoverflow"
$tr1.Font.Size = 20
$tr1.Font.Name = "Arial Black"
$tr1.Font.NameFarEast = "Arial Black"
$tr1.Font.NameComplexScript = "Arial Black"
$tr1.ActionSettings(1).Hyperlink.Address = "https://stackoverflow.com"

# --- "wiki" textbox (links to Wikipedia) ------------------------------------
$tb2 = $s.Shapes.AddTextbox(1, (EMU 9606366), (EMU 6404790), (EMU 768159), (EMU 400110))
$tb2.Name = "TextBox 3"
$tb2.Fill.Visible = 0
$tb2.TextFrame.WordWrap = 0
$tb2.TextFrame.AutoSize = 1

$tr2 = $tb2.TextFrame.TextRange
$tr2.Text = "wiki"
$tr2.Font.Size = 20
$tr2.Font.Name = "Arial Black"
$tr2.Font.NameFarEast = "Arial Black"
$tr2.Font.NameComplexScript = "Arial Black"
$tr2.ActionSettings(1).Hyperlink.Address = "https://en.wikipedia.org/wiki/Multivariate_normal_distribution"
